# Re-pull data: update column F (dSF) values for the scherzer_max.xlsx sheet.
# These are the "final" values corresponding to each start date's dS (delta strikeout)
# stat, refreshed from the source after the game's box score became final.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    4  = -1
    7  = -2
    8  = 4
    9  = 6
    11 = 2
    12 = 4
    13 = 2
    14 = 7
    16 = 5
    19 = 2
    20 = 2
    21 = 4
    24 = 1
    25 = -2
    26 = -2
    27 = 2
    28 = 11
    29 = 1
    30 = 2
    31 = -4
    32 = 1
    34 = -1
    36 = -4
    39 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
